# Generate Report for Handoff
# Update the "Latest Handoff Date" / "Latest Handback DateTime" values for the
# b0d3be6e-e3e7-4550-816c-a5ca6d59ebab row on each sheet to reflect the new
# handback timestamps produced by the report generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-26-09 10:26:37"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-09 10:26:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-09 10:26:37"
